# Updates the "cryptos" price/volume table to the latest scraped values.
# Column D ("Price") values are stored as literal text (they use "."
# as a thousands separator in some rows, e.g. "43.919.42"), so before
# writing any Price cell we force the cell's number format to Text ("@")
# to stop Excel's automatic "looks like a number" conversion from turning
# values such as "241.53" or "102.00" into real numbers (which would also
# silently drop significant trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.919.42"
$ws.Range("E2").Value = "  +1.39%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.359.64"
$ws.Range("E3").Value = "  +0.05%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.34%  "

# Row 5 - XRP
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.691"
$ws.Range("E5").Value = "  +6.12%  "

# Row 6 - BNB
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.53"

# Row 7 - Solana
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "76.43"
$ws.Range("E7").Value = "  +5.30%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.06%  "

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.637"
$ws.Range("E9").Value = "  +27.71%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +4.96%  "

# Row 11 - OKB
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.38"
$ws.Range("E11").Value = "  +0.86%  "

# Row 12 - Avalanche
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "33.43"
$ws.Range("E12").Value = "  +22.41%  "

# Row 13 - Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.57"
$ws.Range("E13").Value = "  +20.17%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +1.86%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.711.30"
$ws.Range("E15").Value = "  +0.10%  "

# Row 16 - Chainlink
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.92"
$ws.Range("E16").Value = "  +4.18%  "

# Row 17 - Polygon
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.923"
$ws.Range("E17").Value = "  +6.75%  "

# Row 18 - WrappedEther
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.351.98"
$ws.Range("E18").Value = "  +0.11%  "

# Row 19 - WrappedBTC
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.878.86"
$ws.Range("E19").Value = "  +1.55%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  +4.61%  "

# Row 21 - Uniswap
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.68"
$ws.Range("E21").Value = "  +4.87%  "

# Row 22 - Litecoin
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "77.51"
$ws.Range("E22").Value = "  +3.18%  "

# Row 23 - BitcoinCash
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "257.04"
$ws.Range("E23").Value = "  +2.43%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  +0.03%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +3.54%  "

# Row 26 - Cosmos
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.21"
$ws.Range("E26").Value = "  +11.52%  "

# Row 27 - WEMIXToken
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.60"
$ws.Range("E27").Value = "  -6.88%  "

# Row 28 - ImmutableX
$ws.Range("E28").Value = "  +14.94%  "

# Row 29 - Toncoin
$ws.Range("E29").Value = "  +1.32%  "

# Row 30 - EthereumClassic
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "23.20"
$ws.Range("E30").Value = "  +2.72%  "

# Row 31 - Monero
$ws.Range("E31").Value = "  +1.16%  "

# Row 32 - Kaspa
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.130"
$ws.Range("E32").Value = "  -2.14%  "

# Row 33 - Stellar
$ws.Range("E33").Value = "  +5.70%  "

# Row 34 - Filecoin
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.32"
$ws.Range("E34").Value = "  +6.18%  "

# Row 35 - Hedera
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0753"
$ws.Range("E35").Value = "  +8.93%  "

# Row 36 - InternetComputer(DFINITY)
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.36"
$ws.Range("E36").Value = "  +6.19%  "

# Row 37 - RenderToken
$ws.Range("E37").Value = "  +2.26%  "

# Row 38 - LidoDAOToken
$ws.Range("E38").Value = "  +0.24%  "

# Row 39 - THORChain
$ws.Range("E39").Value = "  -0.72%  "

# Row 40 - VeChain
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0276"
$ws.Range("E40").Value = "  +8.05%  "

# Row 41 - InjectiveProtocol
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "19.27"
$ws.Range("E41").Value = "  -0.28%  "

# Row 42 - Algorand
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.203"
$ws.Range("E42").Value = "  +18.37%  "

# Row 43 - BinanceUSD
$ws.Range("E43").Value = "  -0.08%  "

# Row 44 - FraxShare
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.92"
$ws.Range("E44").Value = "  -0.37%  "

# Row 45 - Cronos
$ws.Range("E45").Value = "  +6.23%  "

# Row 46 - was TrustWalletToken, now NEARProtocol (rows 46/47 swapped places
# in the ranking and got refreshed values)
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.52"
$ws.Range("E46").Value = "  +13.62%  "

# Row 47 - was NEARProtocol, now TrustWalletToken
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.26"
$ws.Range("E47").Value = "  +5.54%  "

# Row 48 - ARBITRUM
$ws.Range("E48").Value = "  +2.96%  "

# Row 49 - Aave
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "102.00"
$ws.Range("E49").Value = "  +2.54%  "

# Row 50 - FTXToken
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.51"
$ws.Range("E50").Value = "  +0.09%  "

# Row 51 - MultiversX
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.82"
$ws.Range("E51").Value = "  +8.23%  "
